# comment from python script
#
# This script reproduces a data-entry pass across the four worksheets
# (Astronauta, Senador, Mago, Ninja): a new "task" column of partial-credit
# scores (0, 0.5, 0.8, 1, ...) is entered for rows 2-21, one extra grading
# column per sheet (D, D, C, E respectively). On the "Mago" sheet a
# previously-entered score (B10) is also corrected. The running-average
# formulas in the summary column recalculate automatically.

$wb = $excel.ActiveWorkbook

$wsAstronauta = $wb.Worksheets.Item("Astronauta")
$wsSenador    = $wb.Worksheets.Item("Senador")
$wsMago       = $wb.Worksheets.Item("Mago")
$wsNinja      = $wb.Worksheets.Item("Ninja")

# --- Astronauta: new column D (rows 2-21) ---
$astronautaD = @(0, 0.5, 1, 1, 0.8, 0.8, 1, 1, 0, 0, 1, 1, 1, 1, 1, 1, 1, 1, 1, 0.8)
for ($i = 0; $i -lt $astronautaD.Length; $i++) {
    $wsAstronauta.Cells.Item($i + 2, 4).Value = $astronautaD[$i]
}

# --- Senador: new column D (rows 2-21) ---
$senadorD = @(0.5, 0.6, 1, 1, 0.5, 0.5, 0.7, 1, 0, 0, 0.5, 0.6, 1, 1, 1, 1, 0.7, 1, 0.8, 1)
for ($i = 0; $i -lt $senadorD.Length; $i++) {
    $wsSenador.Cells.Item($i + 2, 4).Value = $senadorD[$i]
}

# --- Mago: correct an earlier entry, then new column C (rows 2-21) ---
$wsMago.Cells.Item(10, 2).Value = 0

$magoC = @(0.5, 0.5, 1, 1, 0.5, 0.5, 1, 1, 0, 0, 0.5, 1, 1, 0.8, 1, 1, 1, 1, 1, 1)
for ($i = 0; $i -lt $magoC.Length; $i++) {
    $wsMago.Cells.Item($i + 2, 3).Value = $magoC[$i]
}

# --- Ninja: new column E (rows 2-21) ---
$ninjaE = @(1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 0)
for ($i = 0; $i -lt $ninjaE.Length; $i++) {
    $wsNinja.Cells.Item($i + 2, 5).Value = $ninjaE[$i]
}

# --- Selections / active sheet (Astronauta ends up the active tab) ---
$wsSenador.Range("D20").Select()
$wsMago.Range("C22").Select()
$wsNinja.Range("E4").Select()

$wsAstronauta.Activate()
$wsAstronauta.Range("D6").Select()
